$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextCell $ws "D2" "257.47"
Set-TextCell $ws "E2" "4.99%"
Set-TextCell $ws "D3" "27.67"
Set-TextCell $ws "E3" "-2.67%"
Set-TextCell $ws "D4" "5.231"
Set-TextCell $ws "E4" "-0.16%"
Set-TextCell $ws "D5" "0.05913"
Set-TextCell $ws "E5" "3.74%"
Set-TextCell $ws "D6" "6.691"
Set-TextCell $ws "E6" "1.21%"
Set-TextCell $ws "D7" "0.8697"
Set-TextCell $ws "E7" "2.35%"
Set-TextCell $ws "D8" "1.045"
Set-TextCell $ws "E8" "18.02%"
Set-TextCell $ws "B9" "WazirX"
Set-TextCell $ws "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws "D9" "0.1417"
Set-TextCell $ws "E9" "3.65%"
Set-TextCell $ws "B10" "LiechtensteinCryptoassetsExchange"
Set-TextCell $ws "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws "D10" "0.03637"
Set-TextCell $ws "E10" "9.88%"
Set-TextCell $ws "D11" "0.07200"
Set-TextCell $ws "E11" "2.53%"
Set-TextCell $ws "D12" "0.03260"
Set-TextCell $ws "E12" "3.89%"
Set-TextCell $ws "D13" "0.09221"
Set-TextCell $ws "E13" "0.13%"
Set-TextCell $ws "D14" "0.001543"
Set-TextCell $ws "E14" "0.51%"
Set-TextCell $ws "B15" "One"
Set-TextCell $ws "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws "D15" "0.0006058"
Set-TextCell $ws "E15" "1.54%"
Set-TextCell $ws "B16" "TigerCash"
Set-TextCell $ws "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D16" "0.005887"
Set-TextCell $ws "E16" "-1.35%"
Set-TextCell $ws "B17" "LEO"
Set-TextCell $ws "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D17" "3.482"
Set-TextCell $ws "E17" "-0.30%"
Set-TextCell $ws "B18" "GateToken"
Set-TextCell $ws "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws "D18" "3.266"
Set-TextCell $ws "E18" "2.16%"
Set-TextCell $ws "B19" "BTSEToken"
Set-TextCell $ws "C19" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D19" "2.209"
Set-TextCell $ws "E19" "1.60%"
Set-TextCell $ws "B20" "BitpandaEcosystemToken"
Set-TextCell $ws "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell $ws "D20" "0.3150"
Set-TextCell $ws "E20" "-0.60%"
Set-TextCell $ws "E21" "-0.48%"
Set-TextCell $ws "D22" "3.542"
Set-TextCell $ws "E22" "0.70%"
Set-TextCell $ws "D23" "0.04186"
Set-TextCell $ws "E23" "2.33%"
Set-TextCell $ws "E24" "1.63%"
Set-TextCell $ws "D25" "0.001221"
Set-TextCell $ws "E25" "-0.14%"
Set-TextCell $ws "D26" "0.004534"
Set-TextCell $ws "E26" "9.38%"
Set-TextCell $ws "D27" "0.0001202"
Set-TextCell $ws "E27" "0.19%"
Set-TextCell $ws "D28" "0.0001941"
Set-TextCell $ws "E28" "34.02%"
Set-TextCell $ws "D40" "0.03817"
Set-TextCell $ws "E40" "1.17%"
Set-TextCell $ws "B41" "BKEXToken"
Set-TextCell $ws "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell $ws "D41" "0.1106"
Set-TextCell $ws "E41" "3.79%"
Set-TextCell $ws "B42" "KickToken"
Set-TextCell $ws "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell $ws "D42" "0.003985"
Set-TextCell $ws "E42" "-22.45%"
Set-TextCell $ws "D43" "0.002385"
Set-TextCell $ws "E43" "-4.18%"
Set-TextCell $ws "D44" "0.009913"
Set-TextCell $ws "E44" "7.92%"
Set-TextCell $ws "D45" "0.00005443"
Set-TextCell $ws "E45" "3.25%"
Set-TextCell $ws "E46" "0.20%"
Set-TextCell $ws "E47" "4.11%"
Set-TextCell $ws "E48" "-5.66%"
Set-TextCell $ws "D49" "0.00002103"
Set-TextCell $ws "E49" "0.20%"
Set-TextCell $ws "D50" "0.0002003"
Set-TextCell $ws "E50" "0.20%"
